$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 635
$ws.Range("F5").Value = 810
$ws.Range("F7").Value = 869
$ws.Range("F8").Value = 480
$ws.Range("F9").Value = 7216
$ws.Range("F10").Value = 1928
$ws.Range("F11").Value = 5260
$ws.Range("F13").Value = 307
$ws.Range("F14").Value = 7307
$ws.Range("F15").Value = 8548
$ws.Range("F17").Value = 1125
$ws.Range("F18").Value = 837
$ws.Range("F19").Value = 4310
$ws.Range("F20").Value = 645
$ws.Range("F21").Value = 155
$ws.Range("F22").Value = 79
$ws.Range("F25").Value = 1165
$ws.Range("F26").Value = 75
$ws.Range("F27").Value = 1593
$ws.Range("F28").Value = 675
$ws.Range("F29").Value = 846
$ws.Range("F30").Value = 1825
$ws.Range("F31").Value = 306
$ws.Range("F32").Value = 2193
$ws.Range("F34").Value = 93
$ws.Range("F35").Value = 1397
$ws.Range("F36").Value = 66
$ws.Range("F38").Value = 768
$ws.Range("F39").Value = 376
$ws.Range("F40").Value = 2887
$ws.Range("F41").Value = 3967
$ws.Range("F45").Value = 498
$ws.Range("F47").Value = 852
$ws.Range("F49").Value = 4040

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 54
$ws.Range("F17").Value = 42

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5002

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5002
$ws.Range("F6").Value = 635
$ws.Range("F8").Value = 810
$ws.Range("F10").Value = 869
$ws.Range("F11").Value = 480
$ws.Range("F14").Value = 5260
$ws.Range("F16").Value = 7307
$ws.Range("F17").Value = 7307
$ws.Range("F20").Value = 1125
$ws.Range("F21").Value = 837
$ws.Range("F22").Value = 4310
$ws.Range("F23").Value = 645
$ws.Range("F24").Value = 155
$ws.Range("F28").Value = 1165
$ws.Range("F29").Value = 75
$ws.Range("F30").Value = 1593
$ws.Range("F31").Value = 675
$ws.Range("F32").Value = 846
$ws.Range("F33").Value = 1825
$ws.Range("F34").Value = 306
$ws.Range("F35").Value = 2193
$ws.Range("F39").Value = 768
$ws.Range("F41").Value = 376
$ws.Range("F43").Value = 3967
$ws.Range("F47").Value = 852
$ws.Range("F50").Value = 4040
